# Applies the cryptos-list price/volume refresh described by the commit.
# Columns D ("Price") and E ("Volume(1h)") hold plain text, not numbers, so
# numeric-looking Price strings are written with a leading apostrophe
# (quote-prefix) to force Excel to keep them as text instead of parsing them
# into floating point numbers (which would e.g. drop trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.199.45"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "1.851.43"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'0.6945"
$ws.Range("E5").Value = "  -6.04%  "
$ws.Range("D6").Value = "'238.42"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.3066"
$ws.Range("E8").Value = "  -2.95%  "
$ws.Range("D9").Value = "'0.07570"
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("D10").Value = "'23.48"
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("D11").Value = "'0.08102"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "1.853.01"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "'0.7229"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "'5.182"
$ws.Range("E14").Value = "  -3.93%  "
$ws.Range("D15").Value = "'89.01"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "29.219.62"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "'5.786"
$ws.Range("E17").Value = "  -5.02%  "
$ws.Range("D18").Value = "'240.88"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").Value = "'0.000007713"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").Value = "'13.09"
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "2.095.91"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -5.34%  "
$ws.Range("D25").Value = "'9.014"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "'161.55"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").Value = "'0.1453"
$ws.Range("E27").Value = "  -6.28%  "
$ws.Range("D28").Value = "'18.06"
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("D29").Value = "'1.932"
$ws.Range("E29").Value = "  -4.68%  "
$ws.Range("E30").Value = "  -7.37%  "
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'4.425"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("D33").Value = "'4.040"
$ws.Range("E33").Value = "  -5.16%  "
$ws.Range("D34").Value = "'0.05219"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("D35").Value = "'1.190"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("D36").Value = "'0.7077"
$ws.Range("E36").Value = "  -5.48%  "
$ws.Range("D37").Value = "'0.9997"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'2.664"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").Value = "'0.01860"
$ws.Range("E39").Value = "  -5.36%  "
$ws.Range("D40").Value = "'2.697"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").Value = "'0.9312"
$ws.Range("E41").Value = "  +8.02%  "
$ws.Range("D42").Value = "'5.953"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("D43").Value = "'0.4294"
$ws.Range("E43").Value = "  -5.34%  "
$ws.Range("D44").Value = "1.046.41"
$ws.Range("E44").Value = "  -5.04%  "
$ws.Range("D45").Value = "'69.48"
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").Value = "'7.234"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("D49").Value = "'9.268"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").Value = "'1.736"
$ws.Range("E50").Value = "  -6.55%  "
$ws.Range("D51").Value = "1.991.97"
$ws.Range("E51").Value = "  -2.40%  "
